# Initial update with CPL's work-to-date
# - Re-types the "Year" row labels on the Data sheet (cosmetic whitespace
#   normalisation: a single space before the trailing period instead of
#   several) for every historical row (1968-2019).
# - Adds two new data rows for 2020 and 2021 to the Data sheet, including
#   the "Multiply by to get 2012 Dollars" formula in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Rebuild the "Year" label text used in column A for rows 6 (1968) through
# 57 (2019): a 4-digit year, 77 dots, then " .".
$dots = ""
for ($i = 0; $i -lt 77; $i++) {
    $dots = $dots + "."
}

for ($row = 6; $row -le 57; $row++) {
    $year = 1968 + ($row - 6)
    $label = [string]$year + $dots + " ."
    $ws.Cells.Item($row, 1).Value = $label
}

# New row 58: 2020
$ws.Cells.Item(58, 1).Value = "2020" + $dots + " ."
$ws.Cells.Item(58, 2).Value = 257.557
$ws.Cells.Item(58, 3).Value = 260.065
$ws.Cells.Item(58, 4).Value = 258.811
$ws.Cells.Item(58, 5).Value = 1.4
$ws.Cells.Item(58, 6).Value = 1.2
$ws.Range("G58").Formula = '=$D$50/D58'
$ws.Range("G58").NumberFormat = "0.000"

# New row 59: 2021
$ws.Cells.Item(59, 1).Value = "2021" + $dots + " ."
$ws.Cells.Item(59, 2).Value = 266.236
$ws.Cells.Item(59, 3).Value = 275.703
$ws.Cells.Item(59, 4).Value = 270.97
$ws.Cells.Item(59, 5).Value = 7
$ws.Cells.Item(59, 6).Value = 4.7
$ws.Range("G59").Formula = '=$D$50/D59'
$ws.Range("G59").NumberFormat = "0.000"

# Restore/leave the view state close to how the author left it: cursor
# parked near the newly-added rows on the Data sheet, but the About sheet
# left as the active tab/selection when the workbook was saved.
$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsData = $wb.Worksheets.Item("Data")

$wsData.Activate()
$wsData.Range("H58").Select()

$wsAbout.Activate()
$wsAbout.Range("B6").Select()
